$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb1"
$ws.Range("C2").Value = "Itgb6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 26.144619
$ws.Range("H2").Value = 78.433857
$ws.Range("I2").Value = 0.5211737020083955
$ws.Range("J2").Value = 0.5211737020083955
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.6957970000000001
$ws.Range("N2").Value = 2.087391
$ws.Range("O2").Value = 0.9232770860517062
$ws.Range("P2").Value = 0.9232770860517063
$ws.Range("Q2").Value = 18.191347466343
$ws.Range("R2").Value = 163.722127197087
$ws.Range("S2").Value = 0.4811877369170917
$ws.Range("T2").Value = 0.4811877369170918

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb1"
$ws.Range("C3").Value = "Itgb6"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 26.144619
$ws.Range("H3").Value = 78.433857
$ws.Range("I3").Value = 0.5211737020083955
$ws.Range("J3").Value = 0.5211737020083955
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.05781966666666667
$ws.Range("N3").Value = 0.173459
$ws.Range("O3").Value = 0.07672291394829377
$ws.Range("P3").Value = 0.07672291394829378
$ws.Range("Q3").Value = 1.511673155707
$ws.Range("R3").Value = 13.605058401363
$ws.Range("S3").Value = 0.03998596509130383
$ws.Range("T3").Value = 0.03998596509130384

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tgfb1"
$ws.Range("C4").Value = "Itgb6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.91585
$ws.Range("H4").Value = 53.74755
$ws.Range("I4").Value = 0.3571392594830743
$ws.Range("J4").Value = 0.3571392594830742
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.6957970000000001
$ws.Range("N4").Value = 2.087391
$ws.Range("O4").Value = 0.9232770860517062
$ws.Range("P4").Value = 0.9232770860517063
$ws.Range("Q4").Value = 12.46579468245
$ws.Range("R4").Value = 112.19215214205
$ws.Range("S4").Value = 0.329738494810197
$ws.Range("T4").Value = 0.329738494810197

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tgfb1"
$ws.Range("C5").Value = "Itgb6"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.91585
$ws.Range("H5").Value = 53.74755
$ws.Range("I5").Value = 0.3571392594830743
$ws.Range("J5").Value = 0.3571392594830742
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05781966666666667
$ws.Range("N5").Value = 0.173459
$ws.Range("O5").Value = 0.07672291394829377
$ws.Range("P5").Value = 0.07672291394829378
$ws.Range("Q5").Value = 1.03588847505
$ws.Range("R5").Value = 9.32299627545
$ws.Range("S5").Value = 0.02740076467287727
$ws.Range("T5").Value = 0.02740076467287727

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Tgfb1"
$ws.Range("C6").Value = "Itgb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.104416333333333
$ws.Range("H6").Value = 18.313249
$ws.Range("I6").Value = 0.1216870385085301
$ws.Range("J6").Value = 0.1216870385085301
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6957970000000001
$ws.Range("N6").Value = 2.087391
$ws.Range("O6").Value = 0.9232770860517062
$ws.Range("P6").Value = 0.9232770860517063
$ws.Range("Q6").Value = 4.247434571484334
$ws.Range("R6").Value = 38.226911143359
$ws.Range("S6").Value = 0.1123508543244175
$ws.Range("T6").Value = 0.1123508543244175

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Tgfb1"
$ws.Range("C7").Value = "Itgb6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.104416333333333
$ws.Range("H7").Value = 18.313249
$ws.Range("I7").Value = 0.1216870385085301
$ws.Range("J7").Value = 0.1216870385085301
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.05781966666666667
$ws.Range("N7").Value = 0.173459
$ws.Range("O7").Value = 0.07672291394829377
$ws.Range("P7").Value = 0.07672291394829378
$ws.Range("Q7").Value = 0.3529553175878888
$ws.Range("R7").Value = 3.176597858291
$ws.Range("S7").Value = 0.009336184184112669
$ws.Range("T7").Value = 0.00933618418411267

$ws.Range("A8:T10").EntireRow.Delete()